$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("1:1").Delete()
